$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "257.87"
Set-TextValue "E2" "0.70%"
Set-TextValue "G2" "3"
Set-TextValue "D3" "27.22"
Set-TextValue "E3" "-4.31%"
Set-TextValue "G3" "3"
Set-TextValue "D4" "4.858"
Set-TextValue "E4" "-8.21%"
Set-TextValue "G4" "3"
Set-TextValue "D5" "0.05944"
Set-TextValue "E5" "2.58%"
Set-TextValue "G5" "3"
Set-TextValue "D6" "6.702"
Set-TextValue "E6" "0.11%"
Set-TextValue "G6" "3"
Set-TextValue "D7" "0.8682"
Set-TextValue "E7" "-0.34%"
Set-TextValue "G7" "3"
Set-TextValue "D8" "0.9970"
Set-TextValue "E8" "8.64%"
Set-TextValue "G8" "3"
Set-TextValue "D9" "0.1416"
Set-TextValue "E9" "0.09%"
Set-TextValue "G9" "3"
Set-TextValue "B10" "MandalaExchangeToken"
Set-TextValue "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07195"
Set-TextValue "E10" "0.52%"
Set-TextValue "G10" "3"
Set-TextValue "B11" "BitrueCoin"
Set-TextValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.03148"
Set-TextValue "E11" "0.27%"
Set-TextValue "G11" "3"
Set-TextValue "B12" "BitMartToken"
Set-TextValue "C12" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D12" "0.09254"
Set-TextValue "E12" "-0.03%"
Set-TextValue "G12" "3"
Set-TextValue "B13" "BitForexToken"
Set-TextValue "C13" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D13" "0.001538"
Set-TextValue "E13" "0.02%"
Set-TextValue "G13" "3"
Set-TextValue "B14" "One"
Set-TextValue "C14" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D14" "0.0006091"
Set-TextValue "E14" "1.01%"
Set-TextValue "G14" "3"
Set-TextValue "D15" "0.006014"
Set-TextValue "E15" "2.14%"
Set-TextValue "G15" "3"
Set-TextValue "E16" "-0.26%"
Set-TextValue "G16" "3"
Set-TextValue "D17" "3.261"
Set-TextValue "E17" "1.03%"
Set-TextValue "G17" "3"
Set-TextValue "D18" "2.205"
Set-TextValue "E18" "-2.22%"
Set-TextValue "G18" "3"
Set-TextValue "B19" "BitpandaEcosystemToken"
Set-TextValue "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D19" "0.3147"
Set-TextValue "E19" "0.63%"
Set-TextValue "G19" "3"
Set-TextValue "B20" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C20" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D20" "0.03554"
Set-TextValue "E20" "4.41%"
Set-TextValue "G20" "3"
Set-TextValue "D21" "0.1307"
Set-TextValue "E21" "-0.75%"
Set-TextValue "G21" "3"
Set-TextValue "D22" "3.528"
Set-TextValue "E22" "0.22%"
Set-TextValue "G22" "3"
Set-TextValue "D23" "0.04281"
Set-TextValue "E23" "2.43%"
Set-TextValue "G23" "3"
Set-TextValue "E24" "2.76%"
Set-TextValue "G24" "3"
Set-TextValue "D25" "0.001220"
Set-TextValue "E25" "0.06%"
Set-TextValue "G25" "3"
Set-TextValue "D26" "0.004516"
Set-TextValue "E26" "-9.43%"
Set-TextValue "G26" "3"
Set-TextValue "E27" "0.24%"
Set-TextValue "G27" "3"
Set-TextValue "E28" "-22.79%"
Set-TextValue "G28" "3"
Set-TextValue "G29" "3"
Set-TextValue "G30" "3"
Set-TextValue "G31" "3"
Set-TextValue "G32" "3"
Set-TextValue "G33" "3"
Set-TextValue "G34" "3"
Set-TextValue "G35" "3"
Set-TextValue "G36" "3"
Set-TextValue "G37" "3"
Set-TextValue "G38" "3"
Set-TextValue "G39" "3"
Set-TextValue "D40" "0.03829"
Set-TextValue "E40" "-0.53%"
Set-TextValue "G40" "3"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1104"
Set-TextValue "E41" "0.45%"
Set-TextValue "G41" "3"
Set-TextValue "B42" "KickToken"
Set-TextValue "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.003971"
Set-TextValue "E42" "-31.03%"
Set-TextValue "G42" "3"
Set-TextValue "D43" "0.002310"
Set-TextValue "E43" "-0.25%"
Set-TextValue "G43" "3"
Set-TextValue "E44" "-3.81%"
Set-TextValue "G44" "3"
Set-TextValue "D45" "0.00005493"
Set-TextValue "E45" "4.57%"
Set-TextValue "G45" "3"
Set-TextValue "E46" "0.35%"
Set-TextValue "G46" "3"
Set-TextValue "E47" "28.80%"
Set-TextValue "G47" "3"
Set-TextValue "D48" "0.002159"
Set-TextValue "E48" "-0.64%"
Set-TextValue "G48" "3"
Set-TextValue "E49" "0.35%"
Set-TextValue "G49" "3"
Set-TextValue "E50" "0.35%"
Set-TextValue "G50" "3"
Set-TextValue "G51" "3"
